# Applies the changes described in the commit:
#   "Update tab names in all BOMs, fix bi-color LED naming."
#
# 1) Rename the worksheet tab from "CVMod8 V2" to "BOM".
# 2) Fix the bi-color LED part naming: "LED 3mm Flat Bicolor" -> "LED 3mm Dome Bicolor".
# 3) Move the active selection to C13 (where the LED description lives).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the sheet tab.
$ws.Name = "BOM"

# 2) Fix the LED description text in place (C13 holds the bi-color LED description).
$ws.Range("C13").Value = "LED 3mm Dome Bicolor"

# 3) Update the selection/active cell to C13, matching the saved view state.
$ws.Range("C13").Select()
